# Adds Test Data for Russia, Finland and Hungary markets by cloning the
# existing "Denmark" sheet (same layout/styles) three times, renaming the
# copies, and rewriting their market-name / NGC-code cells + repeater list.

$wb = $excel.ActiveWorkbook

function New-CountrySheet {
    param($SourceName, $NewName, $NgcCode, $MarketName, $DropMzxsdr240)

    $src = $wb.Worksheets.Item($SourceName)
    # Copy goes right after the current last sheet in the workbook.
    $src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $NewName

    # B4 (the NGC code cell) must be written before B2 (the market name
    # cell) so new shared-string entries land in NGC-code-then-market order.
    $new.Range("B4").ClearFormats()
    $new.Range("B4").Value = $NgcCode
    $new.Range("B2").Value = $MarketName

    if ($DropMzxsdr240) {
        # Source (Denmark) carries an extra "MZXSDR240" repeater row that
        # this market's list doesn't have.
        $new.Rows.Item(16).Delete()
    }

    # Rows 3-5 wrap the narrower D column text, producing taller rows.
    $new.Rows.Item(3).RowHeight = 28.8
    $new.Rows.Item(4).RowHeight = 28.8
    $new.Rows.Item(5).RowHeight = 28.8

    # Column widths particular to these new sheets.
    $new.Columns.Item(1).ColumnWidth = 24.28
    $new.Columns.Item(2).ColumnWidth = 15.14
    $new.Columns.Item(4).ColumnWidth = 7.63

    return $new
}

$russia = New-CountrySheet "Denmark" "Russia" "NGC-2929/T2899/" "Russia Market" $true
$russia.Range("A1:D19").Select()

$finland = New-CountrySheet "Denmark" "Finland" "NGC-3130/T2942/" "Finland Market" $false
$finland.Range("A1:D19").Select()

$hungary = New-CountrySheet "Denmark" "Hungary" "NGC-3104/T2991/" "Hungary Market" $true
$hungary.Range("I21").Select()

# Hungary (the last added sheet) ends up the active / selected tab.
$hungary.Activate()
